$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()

$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Select()
